# Preparing to debug 3-betting
# Append a new data row (06/21/21) to the Bankrolls sheet and extend the
# four line charts (Fish / Raymond / Scott / Cedric bankroll-over-time) so
# their category & value series cover the new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bankrolls")

$newRow = 21

# --- New row of data -------------------------------------------------
$ws.Cells.Item($newRow, 1).Value = "06/21/21"           # A21 Date
$ws.Cells.Item($newRow, 1).Style = "Normal"              # match existing A3:A20 (no explicit cell style)

$ws.Cells.Item($newRow, 3).Value = 216.44                 # C21 Fish total net
$ws.Cells.Item($newRow, 4).Value = 277.77                 # D21 Fish bankroll
$ws.Cells.Item($newRow, 5).Value = 61.330000000000013     # E21 Fish own money invested

$ws.Cells.Item($newRow, 8).Value = -123.4                 # H21 Raymond total net
$ws.Cells.Item($newRow, 9).Value = 35.630000000000003     # I21 Raymond bankroll
$ws.Cells.Item($newRow, 10).Value = 159.03                # J21 Raymond own money invested

$ws.Cells.Item($newRow, 13).Value = 461.32000000000011    # M21 Scott total net
$ws.Cells.Item($newRow, 14).Value = 471.36999999999989    # N21 Scott bankroll
$ws.Cells.Item($newRow, 15).Value = 10.050000000000001    # O21 Scott own money invested

$ws.Cells.Item($newRow, 18).Value = -40.979999999999968   # R21 Cedric total net
$ws.Cells.Item($newRow, 19).Value = 66.25                 # S21 Cedric bankroll
$ws.Cells.Item($newRow, 20).Value = 107.23                # T21 Cedric own money invested

# --- Extend the four charts so they pick up row 21 --------------------
$charts = $ws.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
    $chart = $charts.Item($i).Chart
    $series = $chart.SeriesCollection()
    for ($j = 1; $j -le $series.Count; $j++) {
        $s = $series.Item($j)
        $s.XValues = $ws.Range("A2:A21")
        $valueCol = $s.Formula -replace '.*,Bankrolls!\$([A-Z]+)\$2:\$[A-Z]+\$20,.*', '$1'
    }
}
